# Apply the "Changed default property value locations" edit:
#  - Rename headers: D1 Power -> Ability_Power, E1 Range -> Ability_Range
#  - Rename header: H1 Movement -> Can_Move (becomes a boolean flag)
#  - Insert a new column I "Move_Range" that holds the old Movement value
#    (old column I "Wait" shifts right to J, which Insert() does automatically)
#  - For every data row, the old H (Movement) value/text moves into the new
#    I (Move_Range) cell, and H becomes a boolean: TRUE if the unit could
#    move (had a numeric Movement value), FALSE if it was "N/A" (text).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

# Insert a blank column before column I; this pushes the old "Wait" column
# (I) to J, and leaves the old "Movement" column (H) in place.
$ws.Columns("I").Insert()

# Fix up header row text.
$ws.Range("D1").Value = "Ability_Power"
$ws.Range("E1").Value = "Ability_Range"
$ws.Range("H1").Value = "Can_Move"
$ws.Range("I1").Value = "Move_Range"

# Walk every data row and move the old Movement value from H into the new
# Move_Range column I, replacing H with a Can_Move boolean.
for ($r = 2; $r -le $lastRow; $r++) {
    $hCell = $ws.Cells.Item($r, 8)   # column H
    $iCell = $ws.Cells.Item($r, 9)   # column I (newly inserted, blank)

    $hVal = $hCell.Value()

    if ($null -eq $hVal) {
        continue
    }

    $isText = $hVal.GetType().Name -eq "String"

    # Copy the old Movement value/text into the new Move_Range column.
    $iCell.Value = $hVal

    # Movement becomes a simple Can_Move boolean: true unless it was "N/A".
    $hCell.Value = -not $isText
}
